# Course-handout PPTX text tweak:
# Slide 27, body placeholder shape: the run
#   "目标系统构建的过程，本实验没有给出修正方案，而是留给有兴趣的同学自己去解决。"
# gets split into three runs, with the middle segment's wording changed
# from "过程" (process) to "流程" (procedure):
#   "目标系统" | "构建的流程，" | "本实验没有给出修正方案，而是留给有兴趣的同学自己去解决。"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(27)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# The substring "构建的过程，" starts at (1-based) character 102 of the
# shape's full text and is 6 characters long. Re-typing it as
# "构建的流程，" both performs the word change and — because it is
# addressed through Characters() — splits the original single run into
# the three runs shown in the diff (the untouched head "目标系统" and
# tail "本实验……解决。" keep their original run properties; the
# freshly (re)typed middle run is the new one).
$middle = $tr.Characters(102, 6)
$middle.Text = "构建的流程，"
